$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "4) Implementacja zabezpieczen..." (Piotr Dudek section)
#   - "błędnego" -> "nieobsługiwanego"
#   - "błędny rozmiar pliku w formacie *.raw" -> "błędne rozmiary mapy"
#   - drop the two w:proofErr elements that wrapped "raw"
# ---------------------------------------------------------------------------
$finder4 = $d.Content
$found4 = $finder4.Find.Execute(
    "4) Implementacja zabezpieczeń przed wprowadzeniem niepoprawnych danych wejściowych (wybór błędnego formatu pliku, błędny rozmiar pliku w formacie *.raw).",
    $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found4) {
    throw "Could not locate paragraph 4) text to update"
}

# Build a brand-new Range object from the found bounds: InsertXML only
# *replaces* the addressed span when used on a freshly constructed Range.
$target4 = $d.Range($finder4.Start, $finder4.End)

$xml4 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">4) Implementacja zabezpieczeń przed wprowadzeniem niepoprawnych danych wejściowych (wybór </w:t></w:r><w:r><w:t>nieobsługiwanego</w:t></w:r><w:r><w:t xml:space="preserve"> formatu pliku, błędn</w:t></w:r><w:r><w:t>e rozmiary mapy</w:t></w:r><w:r><w:t>).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target4.InsertXML($xml4)

# ---------------------------------------------------------------------------
# Change 2: drop items 10), 11), 12) from Piotr Dudek's section, keep the
# "_GoBack" bookmark attached to the end of item 9), and append one new
# empty paragraph at the very end of the document body.
# ---------------------------------------------------------------------------
$finderStart = $d.Content
$null = $finderStart.Find.Execute("9) Utworzenie drugiej prezentacji.", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$tailStart = $finderStart.Start

$finderEnd = $d.Content
$null = $finderEnd.Find.Execute("12) Podział terenu na warstwy, z odpowiednio ponakładanymi różnymi teksturami, na różnych poziomach (bądź też wymieszanymi).", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$tailEnd = $finderEnd.End

$tailRange = $d.Range($tailStart, $tailEnd)

$xmlTail = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0047124E" w:rsidRDefault="00D51A2B" w:rsidP="00D228C9"><w:r><w:t>9</w:t></w:r><w:r w:rsidR="00D228C9"><w:t>) Utworzenie drugiej prezentacji.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$tailRange.InsertXML($xmlTail)
